$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Table 1 -> sheet1.xml
$ws2 = $wb.Worksheets.Item(2)   # Table 2 -> sheet2.xml
$ws3 = $wb.Worksheets.Item(3)   # Table 3 -> sheet3.xml
$ws4 = $wb.Worksheets.Item(4)   # Table 4 -> sheet4.xml

# ---------------------------------------------------------------------------
# Table 1 (sheet1): append new reaction row 13 (Zn + O -> ZnO)
# ---------------------------------------------------------------------------
$ws1.Range("A13").Value = "R22"
$ws1.Range("B13").Value = "Zn"
$ws1.Range("C13").Value = "O"
$ws1.Range("D13").Value = 1
$ws1.Range("E13").Value = 1
$ws1.Range("F13").Value = "ZnO"
$ws1.Range("G13").Value = 12.02455
$ws1.Range("H13").Value = 33554.102

# ---------------------------------------------------------------------------
# Table 2 (sheet2): new rows 40-41 (K values for Zn / ZnO), new column Q (Zn_liq)
# ---------------------------------------------------------------------------
$ws2.Range("A40").Value = "R39"
$ws2.Range("B40").Value = "Zn"
$ws2.Range("C40").Value = 5.1892
$ws2.Range("D40").Value = -6124.14
$ws2.Range("Q40").Value = 1

# ---------------------------------------------------------------------------
# Table 4 (sheet4): new rows 36-38 (ZnO_g, ZnO_l, Zn_l reactions)
# ---------------------------------------------------------------------------
$ws4.Range("A36").Value = "R34"
$ws4.Range("B36").Value = "ZnO_g"
$ws4.Range("C36").Value = 3.1946
$ws4.Range("D36").Value = -10967.931

$ws4.Range("A37").Value = "R35"
$ws4.Range("E37").Value = "1*Zn, 1*O"

# ---------------------------------------------------------------------------
# Table 3 (sheet3): new rows 42-45 (K41-K44 oxide partition coefficients),
# new column M (ZnO)
# ---------------------------------------------------------------------------
$ws3.Range("A42").Value = "K41"
$ws3.Range("A43").Value = "K42"
$ws3.Range("A44").Value = "K43"
$ws3.Range("A45").Value = "K44"

$ws3.Range("B42").Value = "Zn2SiO4"
$ws3.Range("B43").Value = "ZnTiO3"
$ws3.Range("B44").Value = "Zn2TiO4"
$ws3.Range("B45").Value = "ZnAl2O4"

# ---------------------------------------------------------------------------
# Table 2 (sheet2): header for new column Q + second new row
# ---------------------------------------------------------------------------
$ws2.Range("Q1").Value = "Zn_liq"

$ws2.Range("A41").Value = "R40"
$ws2.Range("B41").Value = "ZnO"
$ws2.Range("C41").Value = 3.1946
$ws2.Range("D41").Value = 10967.931
$ws2.Range("Q41").Value = 1

# ---------------------------------------------------------------------------
# Table 4 (sheet4): finish rows 36-38
# ---------------------------------------------------------------------------
$ws4.Range("A38").Value = "R36"
$ws4.Range("B38").Value = "Zn_l"
$ws4.Range("C38").Value = -5.1892
$ws4.Range("D38").Value = 6124.14
$ws4.Range("E38").Value = "1*Zn"

$ws4.Range("B37").Value = "ZnO_l"
$ws4.Range("C37").Value = -12.02455
$ws4.Range("D37").Value = 33554.102

$ws4.Range("E36").Value = "1*Zn, 0.5*O2"

# ---------------------------------------------------------------------------
# Table 3 (sheet3): finish new rows with numeric data + M column
# ---------------------------------------------------------------------------
$ws3.Range("C42").Value = 0.596
$ws3.Range("D42").Value = 1777.9
$ws3.Range("E42").Value = 1
$ws3.Range("M42").Value = 2

$ws3.Range("C43").Value = 2.793
$ws3.Range("D43").Value = -5625.544
$ws3.Range("J43").Value = 1
$ws3.Range("M43").Value = 1

$ws3.Range("C44").Value = -0.1464
$ws3.Range("D44").Value = 3044.1203
$ws3.Range("J44").Value = 1
$ws3.Range("M44").Value = 2

$ws3.Range("C45").Value = -1.27715
$ws3.Range("D45").Value = 4727.51
$ws3.Range("I45").Value = 1
$ws3.Range("M45").Value = 1

$ws3.Range("M1").Value = "ZnO"

# ---------------------------------------------------------------------------
# Table 4 (sheet4): corrected thermodynamic fit constants for rows 31-33
# (SiO2 K-value match against K_from_Gibbs.py, per commit message)
# ---------------------------------------------------------------------------
$ws4.Range("C31").Value = -15.21
$ws4.Range("D31").Value = 36404
$ws4.Range("C32").Value = -4.5888
$ws4.Range("D32").Value = 14548.7
$ws4.Range("C33").Value = -10.7411
$ws4.Range("D33").Value = 26600.3

# ---------------------------------------------------------------------------
# View state: selections / active sheet (Table 4 becomes the active tab)
# ---------------------------------------------------------------------------
$ws1.Range("B14").Select()

$ws2.Range("Q41").Select()

$ws3.Activate()
$ws3.Range("C46").Select()

$ws4.Activate()
$ws4.Range("C34").Select()
